# Auto-generated Excel COM-interop edit script
# Updates CompStat weekly report: volume/date header + crime-stat table (rows 14-27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: volume number + reporting week dates ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Crime complaint table (rows 14-27) ---
$ws.Range("N14").Value = -83.333333333333
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("N15").Value = -36.363636363636
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 700
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 62.5
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = -3.636363636363
$ws.Range("L16").Value = 12.765957446808
$ws.Range("M16").Value = -14.516129032258
$ws.Range("N16").Value = -87.320574162679
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 18
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = -35.714285714285
$ws.Range("I17").Value = 76
$ws.Range("J17").Value = 107
$ws.Range("K17").Value = -28.971962616822
$ws.Range("L17").Value = 22.58064516129
$ws.Range("M17").Value = 35.714285714285
$ws.Range("N17").Value = -56.571428571428
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -27.272727272727
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 88
$ws.Range("K18").Value = -44.318181818181
$ws.Range("L18").Value = 2.083333333333
$ws.Range("M18").Value = -48.958333333333
$ws.Range("N18").Value = -87.719298245614
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -26.923076923076
$ws.Range("I19").Value = 190
$ws.Range("J19").Value = 332
$ws.Range("K19").Value = -42.771084337349
$ws.Range("L19").Value = 28.378378378378
$ws.Range("M19").Value = 71.171171171171
$ws.Range("N19").Value = 31.03448275862
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 56
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 27.272727272727
$ws.Range("L20").Value = 133.333333333333
$ws.Range("M20").Value = 19.148936170212
$ws.Range("N20").Value = -84.986595174262
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 88
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = -19.266055045871
$ws.Range("I21").Value = 432
$ws.Range("J21").Value = 634
$ws.Range("K21").Value = -31.86119873817
$ws.Range("L21").Value = 28.189910979228
$ws.Range("M21").Value = 13.385826771653
$ws.Range("N21").Value = -71.709233791748
$ws.Range("C14").Copy($ws.Range("C22"))
$ws.Range("C14").Copy($ws.Range("D22"))
$ws.Range("E14").Copy($ws.Range("E22"))
$ws.Range("L22").Value = 57.142857142857
$ws.Range("M22").Value = -26.666666666666
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 64.285714285714
$ws.Range("F24").Value = 116
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = 41.463414634146
$ws.Range("I24").Value = 405
$ws.Range("J24").Value = 438
$ws.Range("K24").Value = -7.534246575342
$ws.Range("L24").Value = 42.105263157894
$ws.Range("M24").Value = 49.446494464944
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -10
$ws.Range("F25").Value = 46
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 43.75
$ws.Range("I25").Value = 159
$ws.Range("J25").Value = 140
$ws.Range("K25").Value = 13.571428571428
$ws.Range("L25").Value = 10.416666666666
$ws.Range("M25").Value = -27.397260273972
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -25
$ws.Range("C27").Value = 3
$ws.Range("C36").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 2
$ws.Range("K36").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 50
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 18
$ws.Range("K27").Value = 33.333333333333
$ws.Range("L27").Value = 0
